$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3172.625
$ws.Range("H73").Value = 3172.625
$ws.Range("H76").Value = 7296
$ws.Range("I76").Value = 5944.5
$ws.Range("J76").Value = 9999
$ws.Range("K76").Value = 5944.5
$ws.Range("L76").Value = 9999
$ws.Range("M76").Value = -5629.5
$ws.Range("N76").Value = -10629
$ws.Range("H79").Value = 7296
$ws.Range("I79").Value = 5944.5
$ws.Range("J79").Value = 9999
$ws.Range("K79").Value = 5944.5
$ws.Range("L79").Value = 9999
$ws.Range("M79").Value = -4852.5
$ws.Range("N79").Value = -12183
$ws.Range("H113").Value = 8607
$ws.Range("I113").Value = 7790
$ws.Range("J113").Value = 9968.666999999999
$ws.Range("K113").Value = 7790
$ws.Range("L113").Value = 9968.666999999999
$ws.Range("M113").Value = -4536
$ws.Range("N113").Value = -16476.667
$ws.Range("H118").Value = 189.5
$ws.Range("J118").Value = 189
$ws.Range("L118").Value = 567
$ws.Range("N118").Value = -3881
$ws.Range("H125").Value = 1437
$ws.Range("I125").Value = 950
$ws.Range("J125").Value = 1924
$ws.Range("K125").Value = 8550
$ws.Range("L125").Value = 17316
$ws.Range("M125").Value = -6090
$ws.Range("N125").Value = -22236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 1998
$ws.Range("I17").Value = 1998
$ws.Range("K17").Value = 1998
$ws.Range("M17").Value = -1825
$ws.Range("H97").Value = 1199
$ws.Range("I97").Value = 1009.8889
$ws.Range("K97").Value = 1009.8889
$ws.Range("M97").Value = -513.8889
$ws.Range("H110").Value = 875.125
$ws.Range("I110").Value = 714.4286
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 714.4286
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1330.5714
$ws.Range("N110").Value = -6090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 123
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 616
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 2009.4286
$ws.Range("I134").Value = 2009.4286
$ws.Range("K134").Value = 6028.2858
$ws.Range("M134").Value = -3493.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1387.2222
$ws.Range("J16").Value = 2166.6667
$ws.Range("L16").Value = 2166.6667
$ws.Range("N16").Value = -2740.6667
$ws.Range("H25").Value = 2600
$ws.Range("I25").Value = 2500
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -2326
$ws.Range("N25").Value = -3348
$ws.Range("H31").Value = 1800
$ws.Range("I31").Value = 1800
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1505
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1800
$ws.Range("I34").Value = 1800
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1598
$ws.Range("N34").ClearContents()
$ws.Range("H86").Value = 5075
$ws.Range("J86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5075
$ws.Range("J89").Value = 6000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232
$ws.Range("H107").Value = 1099.6666
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120
$ws.Range("H113").Value = 1387.2222
$ws.Range("J113").Value = 2166.6667
$ws.Range("L113").Value = 2166.6667
$ws.Range("N113").Value = -6506.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4139.273
$ws.Range("J34").Value = 4942.5557
$ws.Range("L34").Value = 14827.6671
$ws.Range("N34").Value = -14995.6671
$ws.Range("H39").Value = 24800
$ws.Range("J39").Value = 29750
$ws.Range("L39").Value = 89250
$ws.Range("N39").Value = -89838
$ws.Range("H55").Value = 1322.9688
$ws.Range("J55").Value = 1378.0952
$ws.Range("L55").Value = 4134.2856
$ws.Range("N55").Value = -4488.2856
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 8000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8586
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 2204.5
$ws.Range("I122").Value = 1106.1666
$ws.Range("K122").Value = 3318.4998
$ws.Range("M122").Value = -868.4998000000001
$ws.Range("H132").Value = 2469.2856
$ws.Range("I132").Value = 1766.6666
$ws.Range("K132").Value = 5299.9998
$ws.Range("M132").Value = -2769.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7250
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H18").Value = 8000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8344
$ws.Range("H32").Value = 4793.625
$ws.Range("I32").Value = 1907
$ws.Range("K32").Value = 1907
$ws.Range("M32").Value = -1590
$ws.Range("H126").Value = 7250
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 10526
$ws.Range("J64").Value = 10526
$ws.Range("L64").Value = 10526
$ws.Range("N64").Value = -11022
$ws.Range("H67").Value = 10526
$ws.Range("J67").Value = 10526
$ws.Range("L67").Value = 10526
$ws.Range("N67").Value = -12242
